$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A batch of cell-count samples (counted 2023-11-26, initials "LG") had been
# logged with only a sample id in column A, leaving B:J blank and the column L
# CV formula (=STDEV(.)/AVERAGE(.)) evaluating to #DIV/0!. Fill in the
# recorded date, initials, squares-counted, and the six hemocytometer counts
# for each sample so the existing formulas recalculate correctly.
$samples = @(
    @{ Row = 27; SquaresCounted = 1; C1 = 279; C2 = 289; C3 = 315; C4 = 283; C5 = 306; C6 = 326 },
    @{ Row = 35; SquaresCounted = 1; C1 = 210; C2 = 163; C3 = 183; C4 = 188; C5 = 189; C6 = 171 },
    @{ Row = 37; SquaresCounted = 1; C1 = 213; C2 = 171; C3 = 201; C4 = 253; C5 = 226; C6 = 180 },
    @{ Row = 43; SquaresCounted = 1; C1 = 132; C2 = 158; C3 = 135; C4 = 124; C5 = 141; C6 = 134 },
    @{ Row = 44; SquaresCounted = 1; C1 = 259; C2 = 251; C3 = 282; C4 = 258; C5 = 238; C6 = 230 },
    @{ Row = 51; SquaresCounted = 1; C1 = 179; C2 = 169; C3 = 192; C4 = 178; C5 = 213; C6 = 191 },
    @{ Row = 65; SquaresCounted = 2; C1 = 134; C2 = 128; C3 = 128; C4 = 126; C5 = 106; C6 = 124 },
    @{ Row = 84; SquaresCounted = 3; C1 = 184; C2 = 177; C3 = 223; C4 = 220; C5 = 222; C6 = 201 }
)

foreach ($sample in $samples) {
    $r = $sample.Row
    $ws.Range("B$r").Value = 20231126
    $ws.Range("C$r").Value = "LG"
    $ws.Range("D$r").Value = $sample.SquaresCounted
    $ws.Range("E$r").Value = $sample.C1
    $ws.Range("F$r").Value = $sample.C2
    $ws.Range("G$r").Value = $sample.C3
    $ws.Range("H$r").Value = $sample.C4
    $ws.Range("I$r").Value = $sample.C5
    $ws.Range("J$r").Value = $sample.C6
}

# Column L already carries the STDEV/AVERAGE CV formula for every one of the
# rows above (it was simply erroring out on blank inputs); filling the data
# in place lets it recalculate on its own, so it is intentionally left alone
# here -- except for row 35, where the outlier in column E (210, versus the
# ~163-189 range of F:J) is excluded from the CV calculation by pointing the
# formula at F35:J35 instead.
$ws.Range("L35").Formula = "=STDEV(F35:J35)/AVERAGE(F35:J35)"

# The last cell edited interactively was B45; leave the selection there.
$ws.Range("B45").Select()
